$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value2 = '69.692.77'
$ws.Range("E2").Value2 = '  +0.14%  '
$ws.Range("D3").Value2 = '2.512.15'
$ws.Range("E3").Value2 = '  +0.12%  '
$ws.Range("E4").Value2 = '  +0.01%  '
Set-TextValue $ws.Range("D5") '575.74'
$ws.Range("E5").Value2 = '  +0.12%  '
Set-TextValue $ws.Range("D6") '166.98'
$ws.Range("E6").Value2 = '  +0.05%  '
$ws.Range("E7").Value2 = '  -0.05%  '
$ws.Range("E8").Value2 = '  -0.02%  '
$ws.Range("D9").Value2 = '2.510.84'
$ws.Range("E9").Value2 = '  +0.09%  '
$ws.Range("E10").Value2 = '  +1.35%  '
$ws.Range("E11").Value2 = '  -0.33%  '
Set-TextValue $ws.Range("D12") '0.357'
$ws.Range("E12").Value2 = '  +4.34%  '
$ws.Range("E13").Value2 = '  +1.88%  '
$ws.Range("D14").Value2 = '2.971.36'
$ws.Range("E14").Value2 = '  +0.06%  '
$ws.Range("E15").Value2 = '  +1.88%  '
$ws.Range("D16").Value2 = '69.714.66'
$ws.Range("E16").Value2 = '  +0.37%  '
Set-TextValue $ws.Range("D17") '24.93'
$ws.Range("E17").Value2 = '  +0.25%  '
$ws.Range("D18").Value2 = '2.518.56'
$ws.Range("E18").Value2 = '  +0.42%  '
Set-TextValue $ws.Range("D19") '11.30'
$ws.Range("E19").Value2 = '  -1.02%  '
Set-TextValue $ws.Range("D20") '7.54'
$ws.Range("E20").Value2 = '  -2.92%  '
Set-TextValue $ws.Range("D21") '350.59'
$ws.Range("E21").Value2 = '  +0.16%  '
$ws.Range("E22").Value2 = '  -0.62%  '
$ws.Range("E23").Value2 = '  -0.70%  '
$ws.Range("E24").Value2 = '  -0.05%  '
Set-TextValue $ws.Range("D25") '70.26'
$ws.Range("E25").Value2 = '  +2.22%  '
$ws.Range("E26").Value2 = '  -1.21%  '
$ws.Range("E27").Value2 = '  -1.54%  '
$ws.Range("D28").Value2 = '2.636.82'
Set-TextValue $ws.Range("D29") '0.997'
$ws.Range("D30").Replace('0893', '0894') | Out-Null
$ws.Range("E30").Value2 = '  -0.91%  '
$ws.Range("E31").Value2 = '  -0.80%  '
Set-TextValue $ws.Range("D32") '461.99'
$ws.Range("E32").Value2 = '  -3.07%  '
$ws.Range("E33").Value2 = '  -4.42%  '
$ws.Range("E34").Value2 = '  -0.74%  '
$ws.Range("E35").Value2 = '  +0.01%  '
Set-TextValue $ws.Range("D36") '158.89'
$ws.Range("E36").Value2 = '  +2.89%  '
$ws.Range("E37").Value2 = '  +0.28%  '
$ws.Range("E38").Value2 = '  +0.72%  '
Set-TextValue $ws.Range("D39") '18.52'
$ws.Range("E39").Value2 = '  -0.09%  '
$ws.Range("E40").Value2 = '  +0.04%  '
$ws.Range("E41").Value2 = '  +0.02%  '
$ws.Range("E42").Value2 = '  -1.05%  '
$ws.Range("E43").Value2 = '  -0.12%  '
Set-TextValue $ws.Range("D44") '38.16'
$ws.Range("E44").Value2 = '  +0.09%  '
Set-TextValue $ws.Range("D45") '2.23'
$ws.Range("E45").Value2 = '  -4.11%  '
$ws.Range("E46").Value2 = '  -7.96%  '
Set-TextValue $ws.Range("D47") '142.58'
$ws.Range("E47").Value2 = '  -1.06%  '
Set-TextValue $ws.Range("D48") '3.48'
$ws.Range("E48").Value2 = '  -1.65%  '
$ws.Range("E49").Value2 = '  -1.71%  '
$ws.Range("E50").Value2 = '  +0.41%  '
$ws.Range("B51").Value2 = 'THORChain'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D51") '5.82'
$ws.Range("E51").Value2 = '  +4.70%  '
